$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 34, pushing the existing
# rows 34..64 down to 37..67 (values/styles shift automatically, exactly
# as Excel's Rows.Insert behaves).
$ws.Range("A34:A36").EntireRow.Insert()

# Common / repeated field values for the three newly-inserted rows.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/kilo (en caja de 15 kilos)"
$origen    = "Provincia de Limarí"

# Row 34
$r = 34
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44447
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productoId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 100
$ws.Cells.Item($r,14).Value = 3100
$ws.Cells.Item($r,15).Value = 3100
$ws.Cells.Item($r,16).Value = 3100
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 3100
$ws.Cells.Item($r,20).Value = 1

# Row 35
$r = 35
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44447
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productoId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Extra (doble especial)"
$ws.Cells.Item($r,13).Value = 100
$ws.Cells.Item($r,14).Value = 3500
$ws.Cells.Item($r,15).Value = 3500
$ws.Cells.Item($r,16).Value = 3500
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 3500
$ws.Cells.Item($r,20).Value = 1

# Row 36
$r = 36
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44447
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productoId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 330
$ws.Cells.Item($r,14).Value = 2500
$ws.Cells.Item($r,15).Value = 2800
$ws.Cells.Item($r,16).Value = 2659
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 2659
$ws.Cells.Item($r,20).Value = 1
